$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits right after
#    the Wireframe hyperlink paragraph.
# -----------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# -----------------------------------------------------------------------
# 2) Insert the sentence "These two factors make the page's content to
#    look cohesive. " right after "...background. " (before the lone
#    space run that precedes "Instead of a background...").
#    We bracket the freshly inserted text with two temporary bookmarks so
#    that the run gets split cleanly into its own <w:r> element, then we
#    remove the temporary bookmarks again.
# -----------------------------------------------------------------------
$range = $d.Content
$found = $range.Find.Execute("consistent header and footer as well as background. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find insertion anchor for the 'cohesive' sentence"
}
$insertPos = $range.End
$insertText = "These two factors make the page" + [char]0x2019 + "s content to look cohesive. "

$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter($insertText)
$newEndPos = $insertPos + $insertText.Length

$rightAnchor = $d.Range($newEndPos, $newEndPos)
$d.Bookmarks.Add("ZZTmpRight", $rightAnchor)
$leftAnchor = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("ZZTmpLeft", $leftAnchor)

$d.Bookmarks("ZZTmpLeft").Delete()
$d.Bookmarks("ZZTmpRight").Delete()

# -----------------------------------------------------------------------
# 3) Split the navigation-bar run right after "...will have its " and put
#    the "_GoBack" bookmark back at that exact spot (this is what Word
#    does automatically to mark the position of the last edit).
# -----------------------------------------------------------------------
$range2 = $d.Content
$found2 = $range2.Find.Execute("Each button to the other pages will have its ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find split anchor in navigation bar sentence"
}
$splitPos = $range2.End
$splitAnchor = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitAnchor)

# -----------------------------------------------------------------------
# 4) Remove the duplicated "These two factors..." sentence that used to
#    trail "...close proximity to another. " further down in the
#    document (it now lives in its new location from step 2).
# -----------------------------------------------------------------------
$range3 = $d.Content
$oldText = " the content will have a close proximity to another. These two factors make the page" + [char]0x2019 + "s content to look cohesive. "
$newText = " the content will have a close proximity to another. "
$found3 = $range3.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found3) {
    throw "Could not find duplicated 'cohesive' sentence to remove"
}
